# Obituary.docx update:
#  - Turn the three "Web Site" URLs in the Memorials table into real hyperlinks
#  - Resize the Memorials table's two grid columns

$d = $word.ActiveDocument
$table = $d.Tables(1)

function Add-UrlHyperlink($rowIndex) {
    $cell = $table.Cell($rowIndex, 2)
    $para = $cell.Range.Paragraphs(1)
    $r = $para.Range
    # Trim the trailing paragraph mark so Hyperlinks.Add replaces the
    # existing run instead of appending a duplicate copy of the text.
    $r.End = $r.End - 1
    $url = $r.Text
    $d.Hyperlinks.Add($r, $url)
}

# Row 2: AKC Canine Health Foundation -> https://www.akcchf.org/
Add-UrlHyperlink 2
# Row 3: LinkedIn Profile -> https://www.linkedin.com/in/ralphhightower/
Add-UrlHyperlink 3
# Row 4: Photography Portfolio -> https://www.flickr.com/photos/ralphhightower/
Add-UrlHyperlink 4

# Resize the two table columns (values in points; 1 pt = 20 twips)
$table.Columns(1).Width = 3813 / 20
$table.Columns(2).Width = 4106 / 20

Write-Output "done"
